$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: "time_taken", styled like the other headers (e.g. E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data cells F2:F12: time_taken values (stored as text, matching source timestamps)
$times = @(
    "2021-10-05 10:52:52.363749",
    "2021-10-05 10:52:52.363761",
    "2021-10-05 10:52:52.363765",
    "2021-10-05 10:52:52.363768",
    "2021-10-05 10:52:52.363771",
    "2021-10-05 10:52:52.363773",
    "2021-10-05 10:52:52.363776",
    "2021-10-05 10:52:52.363779",
    "2021-10-05 10:52:52.363782",
    "2021-10-05 10:52:52.363785",
    "2021-10-05 10:52:52.363788"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $times[$i]
}
